$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first occurrence of this event's row
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2970
$ws1.Range("F5").Value = 72

# Sheet "全部类型" (All types) - duplicate rows aggregating all sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2970
$ws4.Range("F10").Value = 72
